# Update the "取得日時" (fetched datetime) column on the "ランサーズ" sheet
# for all existing data rows (2-17) to the new timestamp, keeping the
# values as plain text strings (matching the original inlineStr cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-06 18:33:59"

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
